$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the trailing block of paragraphs (originally 9-16):
#    empty, "כך מספר...", empty, "מתחילים את FT...", "מחשבים את W-Q...",
#    "חוזרים ל-2 T...", "מאוד פשוט...", empty.
#    (Work from the back of the document first so earlier indices stay valid.)
# ------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(9)
$pEnd = $d.Paragraphs.Item(16)
$rngDel1 = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rngDel1.Delete()

# ------------------------------------------------------------------
# 2. Delete the "כולכם מכירים..." paragraph and the empty paragraph
#    right after it (originally paragraphs 6-7).
# ------------------------------------------------------------------
$pStart2 = $d.Paragraphs.Item(6)
$pEnd2 = $d.Paragraphs.Item(7)
$rngDel2 = $d.Range($pStart2.Range.Start, $pEnd2.Range.End)
$rngDel2.Delete()

# ------------------------------------------------------------------
# 3. Replace the body text of what is now paragraph 6 (the
#    "כאמור LoRA..." paragraph) with the new Table-GPT summary text,
#    keeping the two pairs of manual line breaks. Reuse a single Range
#    object and collapse-to-end after each insert so everything lands
#    in one run, in order.
# ------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$seg = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$seg.Delete()
$seg.Collapse(1)

$seg.InsertAfter("אתם בטח מודעים ליכולות המטורפות של מודלי שפה אבל הם עדיין מתקשים להסתדר עם דאטה טבלאי. המאמר שנסקור היום ב-#shorthebrewpapereviews מציע שיטת טיוב(fine-tuning) של מודלי שפה שבאה להקנות להם יכולת לעבוד עם טבלאות. ")
$seg.Collapse(0)

$seg.InsertAfter([char]11)
$seg.Collapse(0)

$seg.InsertAfter([char]11)
$seg.Collapse(0)

$seg.InsertAfter("קודם כל בואו נבין למה מודלי שפה בעצם מתקשים לשחזר את הביצועים החזקים שלהם בדאטה טבלאי. הסיבה נעוצה בשינויים בין המאפיינים המהותיים של דאטה טבלאי ושפה טבעית. הטקסט הוא חד כיווני (או משמאל לימין או מימין לשמאל כמו עברית) ולעומת זאת לטבלאות מבנה דו-ממדי. דאטה טקסטואלי לא אינווריאנטי לפרמוטציות לעומת רוב הטבלאות שפרמוטציה של עמודות או של שורות אינה משפיעה על תכונות הטבלה. המחברים מציעים לכייל מודל שפה על המשימות שהן אינהרנטיות לטבלאות שהן כמובן מאוד שונות מהמשימות שאנו רואים בעיבוד שפה טבעית. ")
$seg.Collapse(0)

$seg.InsertAfter([char]11)
$seg.Collapse(0)

$seg.InsertAfter([char]11)
$seg.Collapse(0)

$seg.InsertAfter('למשל אחת המשימות שמודל שפה מכויל עליהן היא זיהוי מקומות בטבלה שבהם יש דאטה חסר. משימה אחרת (טיפה יותר מורכבת) היא לאתר שורות בשתי טבלאות המייצגות את אותו ה"ישות״ (entity). עוד משימות טבלאיות היא השלמה ערכים חסרים בטבלה, הפיכה של שאלה מילולית ששאילתה עבור הטבלה ותמצות של תוכן הטבלה. יש כמעט 20 משימות שונות שעליהן מכיילים מודל שפה והמודל המכויל הנושא שם הלא מפתיע TableGPT מציג ביצועים די טובים.')

# ------------------------------------------------------------------
# 4. Remove the manual line break from paragraph 5 (turn it into an
#    empty paragraph) while preserving paragraph boundaries: insert a
#    fresh empty "Normal" paragraph right before it, then delete the
#    original paragraph (with the break) entirely.
# ------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphBefore()
$brPara = $d.Paragraphs.Item(6)
$brPara.Range.Delete()

# ------------------------------------------------------------------
# 5. Update the Hugging Face paper link.
# ------------------------------------------------------------------
$d.Content.Find.Execute("https://huggingface.co/papers/2310.08659", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://huggingface.co/papers/2310.09263", 2)

# ------------------------------------------------------------------
# 6. Update the arXiv paper link.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Paper: https://arxiv.org/abs/2310.08659v4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Paper: https://arxiv.org/abs/2310.09263v1", 2)

# ------------------------------------------------------------------
# 7. Update the review title.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Review 165: [Short] LoftQ: LoRA-Fine-Tuning-Aware Quantization for Large Language Models", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Review 164: [Short] Table-GPT: Table-tuned GPT for Diverse Table Tasks", 2)

Write-Output "done, paragraph count = $($d.Paragraphs.Count)"
